# Corrected excel sheets for application fix issues
$wb = $excel.ActiveWorkbook

# --- Summary sheet: widen current selection from A7:XFD12 to A7:XFD14 ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()
$wsSummary.Range("A7:XFD14").Select()

# --- Repayment schedule sheet: move selection from B7 to A9:XFD9 ---
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Activate()
$wsSchedule.Range("A9:XFD9").Select()

# --- Transactions sheet: correct the ID values in column A ---
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("A2").Value = 67
$wsTransactions.Range("A3").Value = 66

# Restore the originally active sheet/tab so tabSelected stays on Transactions
$wsTransactions.Activate()
